$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-29 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-30 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("44-39=5", $true, $false, $false, $false, $false, $true, 1, $false, "24+59=83", 2) | Out-Null
$d.Content.Find.Execute("63-53=10", $true, $false, $false, $false, $false, $true, 1, $false, "61+10=71", 2) | Out-Null
$d.Content.Find.Execute("41-2=39", $true, $false, $false, $false, $false, $true, 1, $false, "88-67=21", 2) | Out-Null
$d.Content.Find.Execute("99-39=60", $true, $false, $false, $false, $false, $true, 1, $false, "26-12=14", 2) | Out-Null
$d.Content.Find.Execute("2+18=20", $true, $false, $false, $false, $false, $true, 1, $false, "70+6=76", 2) | Out-Null
$d.Content.Find.Execute("56-24=32", $true, $false, $false, $false, $false, $true, 1, $false, "71-3=68", 2) | Out-Null
$d.Content.Find.Execute("57-46=11", $true, $false, $false, $false, $false, $true, 1, $false, "29-26=3", 2) | Out-Null
$d.Content.Find.Execute("64+31=95", $true, $false, $false, $false, $false, $true, 1, $false, "63-56=7", 2) | Out-Null
$d.Content.Find.Execute("1+28=29", $true, $false, $false, $false, $false, $true, 1, $false, "24+42=66", 2) | Out-Null
$d.Content.Find.Execute("45+49=94", $true, $false, $false, $false, $false, $true, 1, $false, "81-50=31", 2) | Out-Null
$d.Content.Find.Execute("66-3=63", $true, $false, $false, $false, $false, $true, 1, $false, "77+12=89", 2) | Out-Null
$d.Content.Find.Execute("36+39=75", $true, $false, $false, $false, $false, $true, 1, $false, "78-38=40", 2) | Out-Null
$d.Content.Find.Execute("27-19=8", $true, $false, $false, $false, $false, $true, 1, $false, "63+16=79", 2) | Out-Null
$d.Content.Find.Execute("11+72=83", $true, $false, $false, $false, $false, $true, 1, $false, "61+19=80", 2) | Out-Null
$d.Content.Find.Execute("40-14=26", $true, $false, $false, $false, $false, $true, 1, $false, "28-11=17", 2) | Out-Null
$d.Content.Find.Execute("48-13=35", $true, $false, $false, $false, $false, $true, 1, $false, "41+46=87", 2) | Out-Null
$d.Content.Find.Execute("73+18=91", $true, $false, $false, $false, $false, $true, 1, $false, "12+60=72", 2) | Out-Null
$d.Content.Find.Execute("21+15=36", $true, $false, $false, $false, $false, $true, 1, $false, "86-14=72", 2) | Out-Null
$d.Content.Find.Execute("92-70=22", $true, $false, $false, $false, $false, $true, 1, $false, "38+28=66", 2) | Out-Null
$d.Content.Find.Execute("35+7=42", $true, $false, $false, $false, $false, $true, 1, $false, "45-7=38", 2) | Out-Null
$d.Content.Find.Execute("47+28=75", $true, $false, $false, $false, $false, $true, 1, $false, "73+1=74", 2) | Out-Null
$d.Content.Find.Execute("43-29=14", $true, $false, $false, $false, $false, $true, 1, $false, "79-33=46", 2) | Out-Null
$d.Content.Find.Execute("91-49=42", $true, $false, $false, $false, $false, $true, 1, $false, "40-39=1", 2) | Out-Null
$d.Content.Find.Execute("29+68=97", $true, $false, $false, $false, $false, $true, 1, $false, "16+3=19", 2) | Out-Null
$d.Content.Find.Execute("80-49=31", $true, $false, $false, $false, $false, $true, 1, $false, "18-2=16", 2) | Out-Null
$d.Content.Find.Execute("28+27=55", $true, $false, $false, $false, $false, $true, 1, $false, "77+1=78", 2) | Out-Null
$d.Content.Find.Execute("75-54=21", $true, $false, $false, $false, $false, $true, 1, $false, "56+29=85", 2) | Out-Null
$d.Content.Find.Execute("6+91=97", $true, $false, $false, $false, $false, $true, 1, $false, "52-0=52", 2) | Out-Null
$d.Content.Find.Execute("6+64=70", $true, $false, $false, $false, $false, $true, 1, $false, "4+86=90", 2) | Out-Null
$d.Content.Find.Execute("91-54=37", $true, $false, $false, $false, $false, $true, 1, $false, "77-61=16", 2) | Out-Null
$d.Content.Find.Execute("62-1=61", $true, $false, $false, $false, $false, $true, 1, $false, "92-91=1", 2) | Out-Null
$d.Content.Find.Execute("44+39=83", $true, $false, $false, $false, $false, $true, 1, $false, "30+12=42", 2) | Out-Null
$d.Content.Find.Execute("29+69=98", $true, $false, $false, $false, $false, $true, 1, $false, "75-33=42", 2) | Out-Null
$d.Content.Find.Execute("9+65=74", $true, $false, $false, $false, $false, $true, 1, $false, "47-14=33", 2) | Out-Null
$d.Content.Find.Execute("73-59=14", $true, $false, $false, $false, $false, $true, 1, $false, "55+10=65", 2) | Out-Null
$d.Content.Find.Execute("73-10=63", $true, $false, $false, $false, $false, $true, 1, $false, "17-3=14", 2) | Out-Null
$d.Content.Find.Execute("51-50=1", $true, $false, $false, $false, $false, $true, 1, $false, "80-39=41", 2) | Out-Null
$d.Content.Find.Execute("70-21=49", $true, $false, $false, $false, $false, $true, 1, $false, "8+40=48", 2) | Out-Null
$d.Content.Find.Execute("48+34=82", $true, $false, $false, $false, $false, $true, 1, $false, "86-16=70", 2) | Out-Null
$d.Content.Find.Execute("96-93=3", $true, $false, $false, $false, $false, $true, 1, $false, "49-19=30", 2) | Out-Null
$d.Content.Find.Execute("13+41=54", $true, $false, $false, $false, $false, $true, 1, $false, "88-57=31", 2) | Out-Null
$d.Content.Find.Execute("48-28=20", $true, $false, $false, $false, $false, $true, 1, $false, "45+16=61", 2) | Out-Null
$d.Content.Find.Execute("12+2=14", $true, $false, $false, $false, $false, $true, 1, $false, "88-12=76", 2) | Out-Null
$d.Content.Find.Execute("0+98=98", $true, $false, $false, $false, $false, $true, 1, $false, "18+38=56", 2) | Out-Null
$d.Content.Find.Execute("75-22=53", $true, $false, $false, $false, $false, $true, 1, $false, "0+38=38", 2) | Out-Null
$d.Content.Find.Execute("86-74=12", $true, $false, $false, $false, $false, $true, 1, $false, "41+44=85", 2) | Out-Null
$d.Content.Find.Execute("87-27=60", $true, $false, $false, $false, $false, $true, 1, $false, "32+26=58", 2) | Out-Null
$d.Content.Find.Execute("27+3=30", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=91", 2) | Out-Null
$d.Content.Find.Execute("83-43=40", $true, $false, $false, $false, $false, $true, 1, $false, "75+7=82", 2) | Out-Null
$d.Content.Find.Execute("27+29=56", $true, $false, $false, $false, $false, $true, 1, $false, "63+31=94", 2) | Out-Null
$d.Content.Find.Execute("15+21=36", $true, $false, $false, $false, $false, $true, 1, $false, "32+16=48", 2) | Out-Null
$d.Content.Find.Execute("34+2=36", $true, $false, $false, $false, $false, $true, 1, $false, "58-4=54", 2) | Out-Null
$d.Content.Find.Execute("34+55=89", $true, $false, $false, $false, $false, $true, 1, $false, "95-62=33", 2) | Out-Null
$d.Content.Find.Execute("51-48=3", $true, $false, $false, $false, $false, $true, 1, $false, "17+41=58", 2) | Out-Null
$d.Content.Find.Execute("77-14=63", $true, $false, $false, $false, $false, $true, 1, $false, "38+20=58", 2) | Out-Null
$d.Content.Find.Execute("55+28=83", $true, $false, $false, $false, $false, $true, 1, $false, "48-24=24", 2) | Out-Null
$d.Content.Find.Execute("57+32=89", $true, $false, $false, $false, $false, $true, 1, $false, "47-10=37", 2) | Out-Null
$d.Content.Find.Execute("17+25=42", $true, $false, $false, $false, $false, $true, 1, $false, "90-15=75", 2) | Out-Null
$d.Content.Find.Execute("52+1=53", $true, $false, $false, $false, $false, $true, 1, $false, "60+1=61", 2) | Out-Null
$d.Content.Find.Execute("2+38=40", $true, $false, $false, $false, $false, $true, 1, $false, "79+8=87", 2) | Out-Null
$d.Content.Find.Execute("8+32=40", $true, $false, $false, $false, $false, $true, 1, $false, "89-61=28", 2) | Out-Null
$d.Content.Find.Execute("58-26=32", $true, $false, $false, $false, $false, $true, 1, $false, "26+67=93", 2) | Out-Null
$d.Content.Find.Execute("18+69=87", $true, $false, $false, $false, $false, $true, 1, $false, "37-18=19", 2) | Out-Null
$d.Content.Find.Execute("2+45=47", $true, $false, $false, $false, $false, $true, 1, $false, "49+0=49", 2) | Out-Null
$d.Content.Find.Execute("80-34=46", $true, $false, $false, $false, $false, $true, 1, $false, "19+63=82", 2) | Out-Null
$d.Content.Find.Execute("40+49=89", $true, $false, $false, $false, $false, $true, 1, $false, "10+37=47", 2) | Out-Null
$d.Content.Find.Execute("53-41=12", $true, $false, $false, $false, $false, $true, 1, $false, "54+22=76", 2) | Out-Null
$d.Content.Find.Execute("46+49=95", $true, $false, $false, $false, $false, $true, 1, $false, "90-22=68", 2) | Out-Null
$d.Content.Find.Execute("12+56=68", $true, $false, $false, $false, $false, $true, 1, $false, "91+8=99", 2) | Out-Null
$d.Content.Find.Execute("20-11=9", $true, $false, $false, $false, $false, $true, 1, $false, "47-16=31", 2) | Out-Null
$d.Content.Find.Execute("25-17=8", $true, $false, $false, $false, $false, $true, 1, $false, "47+51=98", 2) | Out-Null
$d.Content.Find.Execute("13+32=45", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=82", 2) | Out-Null
$d.Content.Find.Execute("89-41=48", $true, $false, $false, $false, $false, $true, 1, $false, "66-9=57", 2) | Out-Null
$d.Content.Find.Execute("83+15=98", $true, $false, $false, $false, $false, $true, 1, $false, "10+63=73", 2) | Out-Null
$d.Content.Find.Execute("45-19=26", $true, $false, $false, $false, $false, $true, 1, $false, "21-12=9", 2) | Out-Null
$d.Content.Find.Execute("47+36=83", $true, $false, $false, $false, $false, $true, 1, $false, "27+30=57", 2) | Out-Null
$d.Content.Find.Execute("51+4=55", $true, $false, $false, $false, $false, $true, 1, $false, "61-22=39", 2) | Out-Null
$d.Content.Find.Execute("10+47=57", $true, $false, $false, $false, $false, $true, 1, $false, "20+14=34", 2) | Out-Null
$d.Content.Find.Execute("25+23=48", $true, $false, $false, $false, $false, $true, 1, $false, "6+79=85", 2) | Out-Null
$d.Content.Find.Execute("8+76=84", $true, $false, $false, $false, $false, $true, 1, $false, "91-19=72", 2) | Out-Null
$d.Content.Find.Execute("46+6=52", $true, $false, $false, $false, $false, $true, 1, $false, "99-52=47", 2) | Out-Null
$d.Content.Find.Execute("12+36=48", $true, $false, $false, $false, $false, $true, 1, $false, "27+67=94", 2) | Out-Null
$d.Content.Find.Execute("1+54=55", $true, $false, $false, $false, $false, $true, 1, $false, "36+46=82", 2) | Out-Null
$d.Content.Find.Execute("55+26=81", $true, $false, $false, $false, $false, $true, 1, $false, "21+53=74", 2) | Out-Null
$d.Content.Find.Execute("15+54=69", $true, $false, $false, $false, $false, $true, 1, $false, "68-65=3", 2) | Out-Null
$d.Content.Find.Execute("21+3=24", $true, $false, $false, $false, $false, $true, 1, $false, "63-16=47", 2) | Out-Null
$d.Content.Find.Execute("63+34=97", $true, $false, $false, $false, $false, $true, 1, $false, "3+63=66", 2) | Out-Null
$d.Content.Find.Execute("95-11=84", $true, $false, $false, $false, $false, $true, 1, $false, "69-22=47", 2) | Out-Null
$d.Content.Find.Execute("56-16=40", $true, $false, $false, $false, $false, $true, 1, $false, "44-40=4", 2) | Out-Null
$d.Content.Find.Execute("26-8=18", $true, $false, $false, $false, $false, $true, 1, $false, "80+10=90", 2) | Out-Null
$d.Content.Find.Execute("73-39=34", $true, $false, $false, $false, $false, $true, 1, $false, "46-0=46", 2) | Out-Null
$d.Content.Find.Execute("98-63=35", $true, $false, $false, $false, $false, $true, 1, $false, "69+19=88", 2) | Out-Null
$d.Content.Find.Execute("89-51=38", $true, $false, $false, $false, $false, $true, 1, $false, "11+42=53", 2) | Out-Null
$d.Content.Find.Execute("3-2=1", $true, $false, $false, $false, $false, $true, 1, $false, "1+31=32", 2) | Out-Null
$d.Content.Find.Execute("56+9=65", $true, $false, $false, $false, $false, $true, 1, $false, "43-10=33", 2) | Out-Null
$d.Content.Find.Execute("15+34=49", $true, $false, $false, $false, $false, $true, 1, $false, "11+62=73", 2) | Out-Null
$d.Content.Find.Execute("98-62=36", $true, $false, $false, $false, $false, $true, 1, $false, "66+23=89", 2) | Out-Null
$d.Content.Find.Execute("10+40=50", $true, $false, $false, $false, $false, $true, 1, $false, "92-0=92", 2) | Out-Null
$d.Content.Find.Execute("40-13=27", $true, $false, $false, $false, $false, $true, 1, $false, "30-11=19", 2) | Out-Null
$d.Content.Find.Execute("34+40=74", $true, $false, $false, $false, $false, $true, 1, $false, "37+41=78", 2) | Out-Null
